# Updates cryptocurrency price (D) and 1h volume change (E) columns
# for rows 2-51 on Sheet1, per the latest scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "57.781.02"
$ws.Range("E2").Value2 = "  -0.40%  "

$ws.Range("D3").Value2 = "2.456.18"
$ws.Range("E3").Value2 = "  +0.26%  "

$ws.Range("E4").Value2 = "  +0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "510.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -2.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "133.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +2.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  -0.08%  "

$ws.Range("E8").Value2 = "  -1.41%  "

$ws.Range("D9").Value2 = "2.457.51"
$ws.Range("E9").Value2 = "  +0.14%  "

$ws.Range("E10").Value2 = "  -0.17%  "

$ws.Range("E11").Value2 = "  -0.75%  "

$ws.Range("E12").Value2 = "  -0.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "4.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -7.16%  "

$ws.Range("D14").Value2 = "2.892.76"
$ws.Range("E14").Value2 = "  +0.24%  "

$ws.Range("D15").Value2 = "57.743.20"
$ws.Range("E15").Value2 = "  -0.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "21.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +1.18%  "

$ws.Range("E17").Value2 = "  +0.73%  "

$ws.Range("D18").Value2 = "2.434.27"
$ws.Range("E18").Value2 = "  -0.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "10.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -0.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "4.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "315.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  +0.16%  "

$ws.Range("E22").Value2 = "  +4.66%  "

$ws.Range("E23").Value2 = "  +0.12%  "

$ws.Range("E24").Value2 = "  -2.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "65.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +0.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -0.32%  "

$ws.Range("E27").Value2 = "  -1.03%  "

$ws.Range("E28").Value2 = "  -5.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "7.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +4.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "173.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -0.31%  "

$ws.Range("E31").Value2 = "  -0.19%  "

$ws.Range("E32").Value2 = "  -0.04%  "

$ws.Range("E33").Value2 = "  +0.20%  "

$ws.Range("E34").Value2 = "  +0.31%  "

$ws.Range("E35").Value2 = "  +0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "18.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +1.09%  "

$ws.Range("E38").Value2 = "  +5.45%  "

$ws.Range("E39").Value2 = "  +1.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "36.71"
$ws.Range("D40").Style = "Normal"

$ws.Range("E41").Value2 = "  +0.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.809"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -0.26%  "

$ws.Range("E43").Value2 = "  +12.04%  "

$ws.Range("E44").Value2 = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "4.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +2.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "256.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.576"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -1.46%  "

$ws.Range("E48").Value2 = "  -0.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0493"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +0.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.0214"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +1.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "17.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +0.66%  "
